{"js": "// Locate the (single) table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// 1) Simple text replacements on the first four rows.\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\ntable.getCell(3, 0).value = \"161\";\nawait context.sync();\n\n// 2) Insert three brand-new one-cell rows right after row index 3 (the\n//    \"161\" row), matching the three new <w:tr> blocks added by the diff.\ntable.rows.load(\"items\");\nawait context.sync();\ntable.rows.items[3].insertRows(\"After\", 3, [\n  [\"0.00001\"],\n  [\"0.00050\"],\n  [\"0.00015\"],\n]);\nawait context.sync();\n\n// 3) After the insertion, the old rows shift down by 3. Update the ones\n//    whose text changed (old indices 5,6,7,8 -> new indices 8,9,10,11).\ntable.getCell(8, 0).value = \"0.00023\";\ntable.getCell(9, 0).value = \"0.00030\";\ntable.getCell(10, 0).value = \"0.00040\";\ntable.getCell(11, 0).value = \"0.02911\";\nawait context.sync();\n\n// 4) Remove the three rows that immediately follow (old \"0.00018\",\n//    \"0.00021\", \"0.01012\" rows - now at index 12 each time one is removed).\ntable.rows.items[12].delete();\nawait context.sync();\ntable.rows.items[12].delete();\nawait context.sync();\ntable.rows.items[12].delete();\nawait context.sync();\n\n// 5) Collapse the three trailing multi-run/tab-separated rows down to a\n//    single run each, keeping only their first value (net row count is\n//    unchanged since step 2 added 3 and step 4 removed 3).\ntable.getCell(43, 0).value = \"99.96\";\ntable.getCell(44, 0).value = \"0.03\";\ntable.getCell(45, 0).value = \"80\";\nawait context.sync();\n", "ps1": "# The document is already open as $word.ActiveDocument.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1) Simple text replacements on the first four rows (1-based indices).\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"161\"\n\n# 2) Insert three brand-new one-cell rows right after row 4 (the \"161\"\n#    row). Rows.Add(BeforeRow) inserts immediately before the anchor, so\n#    adding repeatedly against the same anchor stacks rows in reverse;\n#    insert them back-to-front to land in forward reading order:\n#    0.00001, 0.00050, 0.00015.\n$anchor = $t.Rows.Item(5)\n$newRow3 = $t.Rows.Add($anchor)\n$newRow3.Cells.Item(1).Range.Text = \"0.00015\"\n$newRow2 = $t.Rows.Add($anchor)\n$newRow2.Cells.Item(1).Range.Text = \"0.00050\"\n$newRow1 = $t.Rows.Add($anchor)\n$newRow1.Cells.Item(1).Range.Text = \"0.00001\"\n\n# 3) The old rows 6,7,8,9 (1-based) shifted down to 9,10,11,12 after\n#    inserting the 3 new rows before old row 5. Update their text.\n$t.Cell(9, 1).Range.Text = \"0.00023\"\n$t.Cell(10, 1).Range.Text = \"0.00030\"\n$t.Cell(11, 1).Range.Text = \"0.00040\"\n$t.Cell(12, 1).Range.Text = \"0.02911\"\n\n# 4) Remove the three rows that used to hold \"0.00018\", \"0.00021\" and\n#    \"0.01012\" - they now all sit at position 13 as each one is deleted.\n$t.Rows.Item(13).Delete()\n$t.Rows.Item(13).Delete()\n$t.Rows.Item(13).Delete()\n\n# 5) Collapse the three trailing multi-run/tab-separated rows down to a\n#    single run each, keeping only their first value. The row count is\n#    unchanged overall (3 inserted in step 2, 3 removed in step 4), so\n#    these are still rows 44, 45, 46 (1-based).\n$t.Cell(44, 1).Range.Text = \"99.96\"\n$t.Cell(45, 1).Range.Text = \"0.03\"\n$t.Cell(46, 1).Range.Text = \"80\"\n"}
